$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Checkout payments")

# --- 1. Capture existing hyperlinks (ref -> target) before we touch anything,
#        because row-insert does NOT auto-shift hyperlink ranges in this engine. ---
$oldCount = $ws.Hyperlinks.Count
$oldRefs = New-Object object[] $oldCount
$oldTargets = New-Object object[] $oldCount
$oldDisplays = New-Object object[] $oldCount
for ($i = 1; $i -le $oldCount; $i++) {
    $hl = $ws.Hyperlinks.Item($i)
    $oldRefs[$i-1] = $hl.Range.Address(0, 0)
    $oldTargets[$i-1] = $hl.Address
    $oldDisplays[$i-1] = $hl.TextToDisplay
}

# Remove all hyperlinks now so re-adding later doesn't collide.
$ws.Hyperlinks.Delete()

# --- 2. Insert the new row 5 (shifts row 5..44 down to 6..45). ---
$ws.Rows.Item(5).Insert()

# --- 3. Populate the new row 5. ---
$ws.Range("A5").Value2 = "New Account Details"
$ws.Range("B5").Value2 = "avayugundla+22@helenoftroy.com"
$ws.Range("C5").Value2 = "avayugundla+22@helenoftroy.com"
$ws.Range("D5").Value2 = "avayugundla+22@helenoftroy.com"
$ws.Range("H5").Value2 = "Lotuswave@123"
$ws.Range("I5").Value2 = "Lotuswave@123"

# --- 4. Re-create all the old hyperlinks, shifting any row >= 5 down by one,
#        since the engine left the stored refs pointing at the pre-insert rows. ---
function Shift-Ref([string]$ref) {
    if ($ref -match '^([A-Z]+)([0-9]+)$') {
        $col = $Matches[1]
        $row = [int]$Matches[2]
        if ($row -ge 5) { $row = $row + 1 }
        return "$col$row"
    }
    return $ref
}

for ($i = 0; $i -lt $oldCount; $i++) {
    $newRef = Shift-Ref $oldRefs[$i]
    $target = $oldTargets[$i]
    $display = $oldDisplays[$i]
    if ([string]::IsNullOrEmpty($display)) {
        $ws.Hyperlinks.Add($ws.Range($newRef), $target) | Out-Null
    } else {
        $ws.Hyperlinks.Add($ws.Range($newRef), $target, "", "", $display) | Out-Null
    }
}

# --- 5. Add the 5 new hyperlinks for the new row. ---
$ws.Hyperlinks.Add($ws.Range("B5"), "mailto:avayugundla+22@helenoftroy.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C5"), "mailto:avayugundla+22@helenoftroy.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D5"), "mailto:avayugundla+22@helenoftroy.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H5"), "mailto:Lotuswave@123") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I5"), "mailto:Lotuswave@123") | Out-Null

# --- 6. Fix up the selection / view (matches the diff: no topLeftCell, new active cell). ---
$ws.Range("I6").Select()

Write-Host "done"
